# "splitsize test adjusted for clarity"
# SplitsizeRes (sheet 1): thin out the splitsize increments from 0.05 steps
# to 0.1 steps, drop the A=0 row, and add a trailing A=1 row (reusing the
# old A=0 row's MSE value). LagData (sheet 2): cosmetic view/style refresh
# only (no data changes) - selection/active-tab moves from LagData back to
# SplitsizeRes.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("SplitsizeRes")
$ws2 = $wb.Worksheets.Item("LagData")

# --- SplitsizeRes: remove the odd (x.x5) rows, from the bottom up so the
# row numbers above each deletion point stay valid --------------------
$ws1.Rows.Item(20).Delete()   # 0.95
$ws1.Rows.Item(18).Delete()   # 0.85
$ws1.Rows.Item(16).Delete()   # 0.75
$ws1.Rows.Item(14).Delete()   # 0.65
$ws1.Rows.Item(12).Delete()   # 0.55
$ws1.Rows.Item(10).Delete()   # 0.45
$ws1.Rows.Item(8).Delete()    # 0.35
$ws1.Rows.Item(6).Delete()    # 0.25
$ws1.Rows.Item(4).Delete()    # 0.15
$ws1.Rows.Item(2).Delete()    # 0 (its MSE value moves to the new last row)

# Remaining rows are now: 0.1, 0.2, 0.3, 0.4, 0.5, 0.6, 0.7, 0.8, 0.9 (rows 2-10).
# Append a new row 11 for Splitsize = 1, reusing row 10's formatting, then
# overwrite the values.
$ws1.Range("A10:B10").Copy()
$ws1.Range("A11:B11").PasteSpecial(-4122)
$ws1.Range("A11").Value = 1
$ws1.Range("B11").Value = 0.37656558733025802

$ws1.Application.CutCopyMode = $false

# Row heights go back to the sheet default (no more explicit 18.75 custom
# height), and the used range shrinks from A1:B20 to A1:B11.
$ws1.Range("A1:B11").Rows.AutoFit()

# Column widths reflow now that the values are simpler (single decimal
# digits instead of the old 0.55000000000000004-style figures).
$ws1.Columns.Item(1).ColumnWidth = 8.21875
$ws1.Columns.Item(2).ColumnWidth = 11.5546875

# --- Styles: the old "general number" numeric format (numFmtId 4, used
# only by the now-deleted A3:A20 style) is no longer needed. Re-format the
# Splitsize column with a single decimal place instead.
$ws1.Range("A2:A11").NumberFormat = "#,##0.0"

# LagData's MSE columns keep the 6-decimal numeric format, just renumbered
# (numFmtId 168 -> 165) after the old numFmtId-4 style above was retired.
$ws2.Range("B2:C5").NumberFormat = "#,##0.000000"

# --- View state: the active sheet/selection moves from LagData back to
# SplitsizeRes. -----------------------------------------------------
$ws2.Range("B39:B40").Select()
$ws1.Activate()
$ws1.Range("A1:B11").Select()
